# river update May 2024
# Remove the obsolete "Chlorophyll A (92nd Percentile)" row (row 7) for the
# 2017 - 2021 period; the rows below it shift up to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Delete()
